# Applies the row-data corrections described by the commit diff.
# The edits only touch specific cells within rows 9/10, 18/19/20 and 26/27 -
# the remaining cells (dates, times, county/municipality names, reporter
# names, etc.) are identical between the rows involved so they are left
# untouched to avoid any unintended type coercion (e.g. Excel turning a
# "2026-02-05" text value into a real date when re-assigned).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-CellValue($row, $col) {
    return $ws.Range($col + $row).Value()
}

function Set-CellValue($row, $col, $value) {
    if ($null -eq $value) {
        $ws.Range($col + $row).Value = ""
    } else {
        $ws.Range($col + $row).Value = $value
    }
}

# ---------------------------------------------------------------------
# Rows 9 <-> 10 : full swap across columns A,B,E,F,G,H,Q,R
# ---------------------------------------------------------------------
$cols_9_10 = @("A","B","E","F","G","H","Q","R")

$row9_old  = @{}
$row10_old = @{}
foreach ($c in $cols_9_10) {
    $row9_old[$c]  = Get-CellValue 9  $c
    $row10_old[$c] = Get-CellValue 10 $c
}
foreach ($c in $cols_9_10) {
    Set-CellValue 9  $c $row10_old[$c]
    Set-CellValue 10 $c $row9_old[$c]
}

# ---------------------------------------------------------------------
# Rows 18, 19, 20 : rotation -> new18 = old20, new19 = old18, new20 = old19
# ---------------------------------------------------------------------
$cols_18_19_20 = @("A","B","E","F","G","H","M","P","Q","R","S","AC","AE")

$row18_old = @{}
$row19_old = @{}
$row20_old = @{}
foreach ($c in $cols_18_19_20) {
    $row18_old[$c] = Get-CellValue 18 $c
    $row19_old[$c] = Get-CellValue 19 $c
    $row20_old[$c] = Get-CellValue 20 $c
}
foreach ($c in $cols_18_19_20) {
    Set-CellValue 18 $c $row20_old[$c]
    Set-CellValue 19 $c $row18_old[$c]
    Set-CellValue 20 $c $row19_old[$c]
}

# ---------------------------------------------------------------------
# Rows 26 <-> 27 : full swap across columns A,Q,R,S,AC
# ---------------------------------------------------------------------
$cols_26_27 = @("A","Q","R","S","AC")

$row26_old = @{}
$row27_old = @{}
foreach ($c in $cols_26_27) {
    $row26_old[$c] = Get-CellValue 26 $c
    $row27_old[$c] = Get-CellValue 27 $c
}
foreach ($c in $cols_26_27) {
    Set-CellValue 26 $c $row27_old[$c]
    Set-CellValue 27 $c $row26_old[$c]
}

Write-Output "Row swaps/rotation applied successfully."
